# Auto-generated Excel COM-interop script applying the Shiva_Profits.xlsx diff.
# Workbook sheets map 1:1 to leve-crafting jobs (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 5 (hunk 0)
$ws.Range("H5").Value = 172
$ws.Range("I5").Value = 200.33333
$ws.Range("K5").Value = 200.33333
$ws.Range("M5").Value = -85.33332999999999

# row 33 (hunk 1)
$ws.Range("H33").Value = 442.3125
$ws.Range("I33").Value = 533.3333
$ws.Range("J33").Value = 169.25
$ws.Range("K33").Value = 533.3333
$ws.Range("L33").Value = 169.25
$ws.Range("M33").Value = -304.3333
$ws.Range("N33").Value = -627.25

# row 137 (hunk 2)
$ws.Range("H137").Value = 5393.3687
$ws.Range("I137").Value = 4966.5
$ws.Range("J137").Value = 5867.6665
$ws.Range("K137").Value = 14899.5
$ws.Range("L137").Value = 17602.9995
$ws.Range("M137").Value = -12349.5
$ws.Range("N137").Value = -22702.9995

# row 138 (hunk 3)
$ws.Range("H138").Value = 5065.5127
$ws.Range("I138").Value = 5072
$ws.Range("J138").Value = 5064.7715
$ws.Range("K138").Value = 15216
$ws.Range("L138").Value = 15194.3145
$ws.Range("M138").Value = -10076
$ws.Range("N138").Value = -25474.3145

$ws = $wb.Worksheets.Item("ARM")
# row 2 (hunk 4)
$ws.Range("H2").Value = 1510.9395
$ws.Range("I2").Value = 1645.5652
$ws.Range("J2").Value = 1201.3
$ws.Range("K2").Value = 1645.5652
$ws.Range("L2").Value = 1201.3
$ws.Range("M2").Value = -1532.5652
$ws.Range("N2").Value = -1427.3

# row 5 (hunk 5)
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

# row 32 (hunk 6)
$ws.Range("H32").Value = 2045.95
$ws.Range("I32").Value = 1307.2809
$ws.Range("K32").Value = 1307.2809
$ws.Range("M32").Value = -1020.2809

# row 74 (hunk 7)
$ws.Range("H74").Value = 1829.3158
$ws.Range("I74").Value = 1423.7858
$ws.Range("J74").Value = 2964.8
$ws.Range("K74").Value = 1423.7858
$ws.Range("L74").Value = 2964.8
$ws.Range("M74").Value = -549.7858000000001
$ws.Range("N74").Value = -4712.8

# row 77 (hunk 8)
$ws.Range("H77").Value = 1829.3158
$ws.Range("I77").Value = 1423.7858
$ws.Range("J77").Value = 2964.8
$ws.Range("K77").Value = 7118.929
$ws.Range("L77").Value = 14824
$ws.Range("M77").Value = -2750.929
$ws.Range("N77").Value = -23560

# row 116 (hunk 9)
$ws.Range("H116").Value = 1510.9395
$ws.Range("I116").Value = 1645.5652
$ws.Range("J116").Value = 1201.3
$ws.Range("K116").Value = 1645.5652
$ws.Range("L116").Value = 1201.3
$ws.Range("M116").Value = 648.4348
$ws.Range("N116").Value = -5789.3

# row 132 (hunk 10)
$ws.Range("H132").Value = 8778.429
$ws.Range("I132").Value = 4788.5
$ws.Range("J132").Value = 18753.25
$ws.Range("K132").Value = 14365.5
$ws.Range("L132").Value = 56259.75
$ws.Range("M132").Value = -11835.5
$ws.Range("N132").Value = -61319.75

$ws = $wb.Worksheets.Item("BSM")
# row 3 (hunk 11)
$ws.Range("H3").Value = 1510.9395
$ws.Range("I3").Value = 1645.5652
$ws.Range("J3").Value = 1201.3
$ws.Range("K3").Value = 1645.5652
$ws.Range("L3").Value = 1201.3
$ws.Range("M3").Value = -1531.5652
$ws.Range("N3").Value = -1429.3

# row 4 (hunk 12)
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

# row 24 (hunk 13)
$ws.Range("H24").Value = 1077
$ws.Range("I24").Value = 1077
$ws.Range("K24").Value = 1077
$ws.Range("M24").Value = -842

# row 29 (hunk 14)
$ws.Range("H29").Value = 1013
$ws.Range("I29").Value = 898.8333
$ws.Range("J29").Value = 1355.5
$ws.Range("K29").Value = 898.8333
$ws.Range("L29").Value = 1355.5
$ws.Range("M29").Value = -609.8333
$ws.Range("N29").Value = -1933.5

# row 33 (hunk 15)
$ws.Range("H33").Value = 2531.7144
$ws.Range("I33").Value = 952.5
$ws.Range("K33").Value = 952.5
$ws.Range("M33").Value = -616.5

# row 47 (hunk 16)
$ws.Range("H47").Value = 110260
$ws.Range("J47").Value = 110260
$ws.Range("L47").Value = 110260
$ws.Range("N47").Value = -111300

# row 94 (hunk 17)
$ws.Range("H94").Value = 6496.154
$ws.Range("I94").Value = 3659.3333
$ws.Range("J94").Value = 8927.714
$ws.Range("K94").Value = 3659.3333
$ws.Range("L94").Value = 8927.714
$ws.Range("M94").Value = -3208.3333
$ws.Range("N94").Value = -9829.714

# row 99 (hunk 18)
$ws.Range("H99").Value = 4281.4585
$ws.Range("I99").Value = 5037.2354
$ws.Range("K99").Value = 5037.2354
$ws.Range("M99").Value = -3539.2354

# row 128 (hunk 19)
$ws.Range("H128").Value = 1989
$ws.Range("I128").Value = 1989
$ws.Range("K128").Value = 5967
$ws.Range("M128").Value = -3477

$ws = $wb.Worksheets.Item("CRP")
# row 31 (hunk 20)
$ws.Range("H31").Value = 3952.7368
$ws.Range("I31").Value = 3585.2727
$ws.Range("J31").Value = 4458
$ws.Range("K31").Value = 3585.2727
$ws.Range("L31").Value = 4458
$ws.Range("M31").Value = -3290.2727
$ws.Range("N31").Value = -5048

# row 34 (hunk 21)
$ws.Range("H34").Value = 3952.7368
$ws.Range("I34").Value = 3585.2727
$ws.Range("J34").Value = 4458
$ws.Range("K34").Value = 3585.2727
$ws.Range("L34").Value = 4458
$ws.Range("M34").Value = -3383.2727
$ws.Range("N34").Value = -4862

# row 58 (hunk 22)
$ws.Range("H58").Value = 6970.8945
$ws.Range("I58").Value = 6937.525
$ws.Range("J58").Value = 7049.4116
$ws.Range("K58").Value = 6937.525
$ws.Range("L58").Value = 7049.4116
$ws.Range("M58").Value = -6734.525
$ws.Range("N58").Value = -7455.4116

# row 132 (hunk 23)
$ws.Range("H132").Value = 2973.075
$ws.Range("I132").Value = 2939.9714
$ws.Range("J132").Value = 3204.8
$ws.Range("K132").Value = 8819.914199999999
$ws.Range("L132").Value = 9614.400000000001
$ws.Range("M132").Value = -6289.914199999999
$ws.Range("N132").Value = -14674.4

# row 136 (hunk 24)
$ws.Range("H136").Value = 6970.8945
$ws.Range("I136").Value = 6937.525
$ws.Range("J136").Value = 7049.4116
$ws.Range("K136").Value = 20812.575
$ws.Range("L136").Value = 21148.2348
$ws.Range("M136").Value = -18262.575
$ws.Range("N136").Value = -26248.2348

$ws = $wb.Worksheets.Item("CUL")
# row 86 (hunk 25)
$ws.Range("H86").Value = 2391.2307
$ws.Range("I86").Value = 383.42856
$ws.Range("K86").Value = 1150.28568
$ws.Range("M86").Value = 35.71432000000004

# row 89 (hunk 26)
$ws.Range("H89").Value = 2391.2307
$ws.Range("I89").Value = 383.42856
$ws.Range("K89").Value = 3450.85704
$ws.Range("M89").Value = 2477.14296

$ws = $wb.Worksheets.Item("GSM")
# row 39 (hunk 27)
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

# row 97 (hunk 28)
$ws.Range("H97").Value = 473
$ws.Range("I97").Value = 473
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 473
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 23
$ws.Range("N97").ClearContents()

# row 132 (hunk 29)
$ws.Range("H132").Value = 6523.921
$ws.Range("I132").Value = 6413.6333
$ws.Range("J132").Value = 6937.5
$ws.Range("K132").Value = 19240.8999
$ws.Range("L132").Value = 20812.5
$ws.Range("M132").Value = -16710.8999
$ws.Range("N132").Value = -25872.5

$ws = $wb.Worksheets.Item("LTW")
# row 6 (hunk 30)
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()

# row 61 (hunk 31)
$ws.Range("H61").Value = 8779773
$ws.Range("I61").Value = 9014.483
$ws.Range("J61").Value = 47621704
$ws.Range("K61").Value = 9014.483
$ws.Range("L61").Value = 47621704
$ws.Range("M61").Value = -8812.483
$ws.Range("N61").Value = -47622108

# row 93 (hunk 32)
$ws.Range("H93").Value = 3740.3333
$ws.Range("I93").Value = 4350.3823
$ws.Range("J93").Value = 2520.2354
$ws.Range("K93").Value = 4350.3823
$ws.Range("L93").Value = 2520.2354
$ws.Range("M93").Value = -3102.3823
$ws.Range("N93").Value = -5016.2354

# row 113 (hunk 33)
$ws.Range("H113").Value = 8779773
$ws.Range("I113").Value = 9014.483
$ws.Range("J113").Value = 47621704
$ws.Range("K113").Value = 9014.483
$ws.Range("L113").Value = 47621704
$ws.Range("M113").Value = -6844.483
$ws.Range("N113").Value = -47626044

# row 132 (hunk 34)
$ws.Range("H132").Value = 3619.1875
$ws.Range("I132").Value = 2833.6667
$ws.Range("J132").Value = 4629.143
$ws.Range("K132").Value = 8501.000100000001
$ws.Range("L132").Value = 13887.429
$ws.Range("M132").Value = -5971.000100000001
$ws.Range("N132").Value = -18947.429

# row 136 (hunk 35)
$ws.Range("H136").Value = 6919.1816
$ws.Range("I136").Value = 2275.8667
$ws.Range("J136").Value = 16869.143
$ws.Range("K136").Value = 6827.6001
$ws.Range("L136").Value = 50607.429
$ws.Range("M136").Value = -4277.6001
$ws.Range("N136").Value = -55707.429

$ws = $wb.Worksheets.Item("WVR")
# row 7 (hunk 36)
$ws.Range("H7").Value = 6590.3335
$ws.Range("I7").Value = 6380
$ws.Range("J7").Value = 6800.6665
$ws.Range("K7").Value = 6380
$ws.Range("L7").Value = 6800.6665
$ws.Range("N7").Value = -7026.6665
$ws.Range("M7").Value = -6267

# row 96 (hunk 37)
$ws.Range("H96").Value = 2664.6667
$ws.Range("I96").Value = 1490
$ws.Range("K96").Value = 1490
$ws.Range("M96").Value = -117

